# Apply the "Added mapping to CV terms based on Goslin shorthand level" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet.
$ws.Name = "cv_term_map"

# Insert a new column before column A to hold the Goslin shorthand-level
# constant (mapped to the existing human readable level name).
$ws.Columns.Item(1).Insert()

$levels = @(
    "CATEGORY",
    "CLASS",
    "SPECIES",
    "PHOSPHATE_POSITION",
    "MOLECULAR_SPECIES",
    "SN_POSITION",
    "DBE_POSITION",
    "STRUCTURE_DEFINED",
    "FULL_STRUCTURE",
    "COMPLETE_STRUCTURE"
)

$ws.Cells.Item(1, 1).Value = "Shorthand.Level"

for ($i = 0; $i -lt $levels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $levels[$i]
}

# Match the widened column width from the diff (both columns share the
# same custom width, same as before the edit just re-applied to column B).
$ws.Columns.Item(1).ColumnWidth = 37.5
$ws.Columns.Item(2).ColumnWidth = 37.5

# Match the new active selection.
$ws.Range("B1").Select()
